# Farrowing.xlsx - StatusView update
# - Re-points the "Boar Used" entries for rows 16-21 to boar 0680
# - Fills in the weaning data (Mortality..Remarks) for rows 66-69
# - Adds four new farrowing records (rows 89-92)
# - Updates the saved window view (top-left cell / active selection)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Farrowing")
$ws.Activate()

# ---------------------------------------------------------------------------
# 1. Boar Used correction for rows 16, 17, 18, 20, 21 -> "0680"
# ---------------------------------------------------------------------------
$ws.Range("E16").Value = "0680"
$ws.Range("E17").Value = "0680"
$ws.Range("E18").Value = "0680"
$ws.Range("E20").Value = "0680"
$ws.Range("E21").Value = "0680"

# Match the plain text style (no border) used elsewhere in the Boar Used column
$ws.Range("C66").Copy() | Out-Null
$ws.Range("E16").PasteSpecial(-4122) | Out-Null
$ws.Range("E17").PasteSpecial(-4122) | Out-Null
$ws.Range("E18").PasteSpecial(-4122) | Out-Null
$ws.Range("E20").PasteSpecial(-4122) | Out-Null
$ws.Range("E21").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# ---------------------------------------------------------------------------
# 2. Weaning info (Mortality, Weaning Date, Weaning Quantity, Weaning Age,
#    Average Weaning Weight, Remarks) for rows 66-69
# ---------------------------------------------------------------------------
$ws.Range("M66").Value = 0
$ws.Range("N66").Value = 43182
$ws.Range("O66").Value = 7
$ws.Range("P66").Value = 30
$ws.Range("Q66").Value = 9.01
$ws.Range("R66").Value = "Weaned"

$ws.Range("M67").Value = 0
$ws.Range("N67").Value = 43182
$ws.Range("O67").Value = 4
$ws.Range("P67").Value = 30
$ws.Range("Q67").Value = 8.62
$ws.Range("R67").Value = "Weaned -4 Heads Adopt(03-23-18)"

$ws.Range("M68").Value = 0
$ws.Range("N68").Value = 43182
$ws.Range("O68").Value = 7
$ws.Range("P68").Value = 30
$ws.Range("Q68").Value = 8.42
$ws.Range("R68").Value = "Weaned +1 Heads Adopt(02-25-18)"

$ws.Range("M69").Value = 0
$ws.Range("N69").Value = 43182
$ws.Range("O69").Value = 3
$ws.Range("P69").Value = 28
$ws.Range("Q69").Value = 7.8
$ws.Range("R69").Value = "Weaned"

# Weaning Date columns use the mm-dd-yy date format already used in N2
$ws.Range("N2").Copy() | Out-Null
$ws.Range("N66").PasteSpecial(-4122) | Out-Null
$ws.Range("N67").PasteSpecial(-4122) | Out-Null
$ws.Range("N68").PasteSpecial(-4122) | Out-Null
$ws.Range("N69").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# Q69 keeps the 0.00 numeric style seen elsewhere in that column
$ws.Range("Q22").Copy() | Out-Null
$ws.Range("Q69").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# ---------------------------------------------------------------------------
# 3. Four new farrowing records - rows 89, 90, 91, 92
# ---------------------------------------------------------------------------
$ws.Range("A89").Value = "A135"
$ws.Range("B89").Value = 43178
$ws.Range("C89").Value = "08570"
$ws.Range("D89").Value = "TP"
$ws.Range("E89").Value = "0678/0694"
$ws.Range("F89").Value = 13
$ws.Range("G89").Value = 11
$ws.Range("H89").Value = 7
$ws.Range("I89").Value = 4
$ws.Range("J89").Value = 2
$ws.Range("K89").Value = 0
$ws.Range("L89").Value = 1.45

$ws.Range("A90").Value = "A137"
$ws.Range("B90").Value = 43179
$ws.Range("C90").Value = "08575"
$ws.Range("D90").Value = "TP"
$ws.Range("E90").Value = "0678"
$ws.Range("F90").Value = 15
$ws.Range("G90").Value = 12
$ws.Range("H90").Value = 6
$ws.Range("I90").Value = 6
$ws.Range("J90").Value = 3
$ws.Range("K90").Value = 0
$ws.Range("L90").Value = 1.35

$ws.Range("A91").Value = "A138"
$ws.Range("B91").Value = 43181
$ws.Range("C91").Value = "09411"
$ws.Range("D91").Value = "TP"
$ws.Range("E91").Value = "0678"
$ws.Range("F91").Value = 11
$ws.Range("G91").Value = 10
$ws.Range("H91").Value = 7
$ws.Range("I91").Value = 3
$ws.Range("J91").Value = 1
$ws.Range("K91").Value = 0
$ws.Range("L91").Value = 1.31

$ws.Range("A92").Value = "A139"
$ws.Range("B92").Value = 43181
$ws.Range("C92").Value = "08982"
$ws.Range("D92").Value = "TP"
$ws.Range("E92").Value = "0678"
$ws.Range("F92").Value = 11
$ws.Range("G92").Value = 9
$ws.Range("H92").Value = 4
$ws.Range("I92").Value = 5
$ws.Range("J92").Value = 1
$ws.Range("K92").Value = 0
$ws.Range("L92").Value = 1.63

# Sow No (fill-highlighted text) style, matching the rest of the table
$ws.Range("C75").Copy() | Out-Null
$ws.Range("C89").PasteSpecial(-4122) | Out-Null
$ws.Range("C90").PasteSpecial(-4122) | Out-Null
$ws.Range("C91").PasteSpecial(-4122) | Out-Null
$ws.Range("C92").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# Boar Used text style (no border)
$ws.Range("C66").Copy() | Out-Null
$ws.Range("E89").PasteSpecial(-4122) | Out-Null
$ws.Range("E90").PasteSpecial(-4122) | Out-Null
$ws.Range("E91").PasteSpecial(-4122) | Out-Null
$ws.Range("E92").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# Litter count columns (F:K) numeric style
$ws.Range("F66:K66").Copy() | Out-Null
$ws.Range("F89:K89").PasteSpecial(-4122) | Out-Null
$ws.Range("F90:K90").PasteSpecial(-4122) | Out-Null
$ws.Range("F91:K91").PasteSpecial(-4122) | Out-Null
$ws.Range("F92:K92").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# ---------------------------------------------------------------------------
# 4. Refresh the saved view - scrolled down with F94 as the active cell
# ---------------------------------------------------------------------------
$excel.ActiveWindow.ScrollRow = 25
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("F94").Select() | Out-Null
